# Edit script: recompute LR-pairs cross table for Gpc3-Igf1r (4 clusters x 4 clusters)
# following the updated clustering/advice ("Natmi following Dr Hou advice").
# Adds the "M2" cluster into both the Sending/Target cluster dimensions, producing a
# full 4x4 (ECs, FAPs, M2, sCs) grid of 16 data rows (rows 2-17) instead of the
# original 3x4 (ECs, FAPs, M2) x (ECs, FAPs, sCs, M2) = 12 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: rowNumber, A(Sending cluster), D(Target cluster),
#             E,F,G,H,I,J,K,L,M,N,O,P,Q,R,S,T
$rowData = @(
    @(2,"ECs","ECs",3,1,1.465847333333333,4.397542,0.0219775736133859,0.02197757361338591,3,1,14.561928,43.685784,0.3501777048818433,0.3501777048818433,21.345563326992,192.110069942928,0.007696056286807236,0.007696056286807237),
    @(3,"ECs","FAPs",3,1,1.465847333333333,4.397542,0.0219775736133859,0.02197757361338591,3,1,14.40015733333333,43.200472,0.3462875276490937,0.3462875276490937,21.10843222664711,189.975890039824,0.007610559630305364,0.007610559630305365),
    @(4,"ECs","M2",3,1,1.465847333333333,4.397542,0.0219775736133859,0.02197757361338591,3,1,2.886742333333333,8.660226999999999,0.06941888497676431,0.06941888497676431,4.23152355133711,38.08371196203399,0.001525658654736007,0.001525658654736007),
    @(5,"ECs","sCs",3,1,1.465847333333333,4.397542,0.0219775736133859,0.02197757361338591,3,1,9.735567333333334,29.206702,0.2341158824922987,0.2341158824922987,14.27085541405378,128.437698726484,0.005145299041537299,0.0051452990415373),
    @(6,"FAPs","ECs",3,1,51.92481233333333,155.774437,0.7785131230699432,0.7785131230699435,3,1,14.561928,43.685784,0.3501777048818433,0.3501777048818433,756.1253786115119,6805.128407503606,0.2726179386570287,0.2726179386570288),
    @(7,"FAPs","FAPs",3,1,51.92481233333333,155.774437,0.7785131230699432,0.7785131230699435,3,1,14.40015733333333,43.200472,0.3462875276490937,0.3462875276490937,747.7254671038071,6729.529203934264,0.2695893846302653,0.2695893846302653),
    @(8,"FAPs","M2",3,1,51.92481233333333,155.774437,0.7785131230699432,0.7785131230699435,3,1,2.886742333333333,8.660226999999999,0.06941888497676431,0.06941888497676431,149.8935539130221,1349.041985217199,0.05404351294329394,0.05404351294329396),
    @(9,"FAPs","sCs",3,1,51.92481233333333,155.774437,0.7785131230699432,0.7785131230699435,3,1,9.735567333333334,29.206702,0.2341158824922987,0.2341158824922987,505.5175067418638,4549.657560676774,0.1822622868393553,0.1822622868393554),
    @(10,"M2","ECs",2,0.6666666666666666,0.05191500000000001,0.155745,0.0007783660059225787,0.0007783660059225788,3,1,14.561928,43.685784,0.3501777048818433,0.3501777048818433,0.7559824921200001,6.80384242908,0.0002725664215120158,0.0002725664215120159),
    @(11,"M2","FAPs",2,0.6666666666666666,0.05191500000000001,0.155745,0.0007783660059225787,0.0007783660059225788,3,1,14.40015733333333,43.200472,0.3462875276490937,0.3462875276490937,0.7475841679600002,6.728257511640002,0.0002695384397970296,0.0002695384397970296),
    @(12,"M2","M2",2,0.6666666666666666,0.05191500000000001,0.155745,0.0007783660059225787,0.0007783660059225788,3,1,2.886742333333333,8.660226999999999,0.06941888497676431,0.06941888497676431,0.149865228235,1.348787054115,0.00005403330023496294,0.00005403330023496294),
    @(13,"M2","sCs",2,0.6666666666666666,0.05191500000000001,0.155745,0.0007783660059225787,0.0007783660059225788,3,1,9.735567333333334,29.206702,0.2341158824922987,0.2341158824922987,0.5054219781100001,4.54879780299,0.0001822278443785703,0.0001822278443785703),
    @(14,"sCs","ECs",3,1,13.25483966666666,39.76451899999999,0.1987309373107482,0.1987309373107483,3,1,14.561928,43.685784,0.3501777048818433,0.3501777048818433,193.016020877544,1737.144187897896,0.06959114351649529,0.0695911435164953),
    @(15,"sCs","FAPs",3,1,13.25483966666666,39.76451899999999,0.1987309373107482,0.1987309373107483,3,1,14.40015733333333,43.200472,0.3462875276490937,0.3462875276490937,190.8717766281075,1717.845989652968,0.06881804494872602,0.06881804494872605),
    @(16,"sCs","M2",3,1,13.25483966666666,39.76451899999999,0.1987309373107482,0.1987309373107483,3,1,2.886742333333333,8.660226999999999,0.06941888497676431,0.06941888497676431,38.26330678731254,344.3697610858129,0.01379568007849939,0.01379568007849939),
    @(17,"sCs","sCs",3,1,13.25483966666666,39.76451899999999,0.1987309373107482,0.1987309373107483,3,1,9.735567333333334,29.206702,0.2341158824922987,0.2341158824922987,129.0433840673709,1161.390456606338,0.04652606876702751,0.04652606876702751)
)

foreach ($row in $rowData) {
    $r = $row[0]
    $ws.Cells.Item($r,1).Value = $row[1]   # A - Sending cluster
    $ws.Cells.Item($r,2).Value = "Gpc3"    # B - Ligand symbol (constant)
    $ws.Cells.Item($r,3).Value = "Igf1r"   # C - Receptor symbol (constant)
    $ws.Cells.Item($r,4).Value = $row[2]   # D - Target cluster

    for ($i = 0; $i -lt 16; $i++) {
        $col = 5 + $i      # columns E..T
        $ws.Cells.Item($r, $col).Value = $row[3 + $i]
    }
}
